# Tissue sheet: reorder columns B,C,D,E (header row 1 + data rows 2..11):
#   new B (user)                     <- old D
#   new C (tissue_name)              <- old B
#   new D (correction_factor_tissue) <- old E
#   new E (approved)                 <- old C
#
# Strategy:
#  1. Stage the original B..E values into a scratch area (H..K) with
#     PasteSpecial-values so text/number typing (e.g. the "TRUE" shared
#     string) is preserved exactly instead of being reinterpreted by .Value
#     (which would turn "TRUE" into a native boolean).
#  2. Use a full-column Cut + "Insert Cut Cells" to relocate column B (the
#     one carrying the custom 24.6328125 width) over to column C. This is
#     how the custom width actually travels with the column in Excel,
#     giving an exact width match instead of the rounding drift you get by
#     re-deriving ColumnWidth (character units) from scratch.
#  3. Overwrite B..E with the staged values in their final target order -
#     this fixes up whatever the column-shift in step 2 left behind, so the
#     end state is correct regardless of how step 2 rearranged things.
#  4. Clean the scratch area and restore the selection/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlShiftToRight = -4161

$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
if (-not $lastRow -or $lastRow -lt 1) { $lastRow = 11 }

# --- 1. Stage original values -------------------------------------------
$ws.Range("B1:E$lastRow").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = 0
# H=old B (tissue_name), I=old C (approved), J=old D (user), K=old E (correction_factor_tissue)

# --- 2. Carry column B's custom width over to column C ------------------
$ws.Columns.Item(2).Cut() | Out-Null
$ws.Columns.Item(4).Insert($xlShiftToRight) | Out-Null
$excel.CutCopyMode = 0

# --- 3. Write the final values back in the correct target order ---------
# new B <- old D (staged at J)
$ws.Range("J1:J$lastRow").Copy() | Out-Null
$ws.Range("B1").PasteSpecial($xlPasteValues) | Out-Null

# new C <- old B (staged at H)
$ws.Range("H1:H$lastRow").Copy() | Out-Null
$ws.Range("C1").PasteSpecial($xlPasteValues) | Out-Null

# new D <- old E (staged at K)
$ws.Range("K1:K$lastRow").Copy() | Out-Null
$ws.Range("D1").PasteSpecial($xlPasteValues) | Out-Null

# new E <- old C (staged at I)
$ws.Range("I1:I$lastRow").Copy() | Out-Null
$ws.Range("E1").PasteSpecial($xlPasteValues) | Out-Null

$excel.CutCopyMode = 0

# --- 4. Cleanup -----------------------------------------------------------
$ws.Range("H1:K$lastRow").Clear() | Out-Null

# Update selection to match the saved state: the whole column D selected,
# active cell D1.
$ws.Range("D1:D1048576").Select()
